# Word COM-interop script implementing the target edit.
#
# Strategy: use Range.InsertXML(...) with exact OOXML fragments for
# paragraph-content changes (gives byte-for-byte control over run/proofErr
# structure), Range.InsertParagraphAfter() to mint new, correctly-styled
# paragraphs (it clones pPr/rPr from the anchor paragraph), and
# Range.Delete() to remove whole paragraphs (the paragraph's Range
# includes its own pilcrow, so deleting it removes the whole paragraph).
#
# Ordering note: paragraph indices renumber as soon as paragraphs are
# deleted/inserted earlier in the body, so every edit below that is
# addressed by absolute Paragraphs(N) index is done against the
# *original* numbering first; the two early empty-paragraph deletions
# (which would shift every later index down by one each) are performed
# last of all.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) & 2) In the table (Process Engineer row), merge the split "Set "/
#    "Test Applications" and "Set "/"Parameters" runs into single runs.
# ---------------------------------------------------------------------
$tbl = $d.Tables(1)
$cell = $tbl.Cell(2, 3)

$p = $cell.Range.Paragraphs(2)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>Set Test Applications</w:t></w:r><w:r><w:t xml:space=`"preserve`"> with order</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

$p = $cell.Range.Paragraphs(3)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>Set Parameters</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) Wrap "prev" with spellcheck proofErr markers in the FAIL/return cell.
# ---------------------------------------------------------------------
$cell = $tbl.Cell(4, 6)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:r><w:t xml:space=`"preserve`">If at least 1 FAIL, return to </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>prev</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> step </w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) Split "Future:" into bold+underlined "Future" plus plain ":".
#    (original paragraph index 46)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(46)
$xml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val=`"single`"/><w:lang w:val=`"en-US`"/></w:rPr><w:t>Future</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>:</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 5) "User authentication" -> "User authentication & roles" (index 47)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(47)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>User authentication</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`"> &amp; roles</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) "Execution order" -> "Versions (*3)" (index 48)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(48)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Versions (*3)</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 7) "Inheritance of item types (based on)" + "   " -> "Deployment to
#    cloud" (index 49)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(49)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Deployment to cloud</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 8) "Export to pdf" -> "Automatic build" (index 50)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(50)
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Automatic build</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 9) Insert 7 new "Future:" bullet paragraphs between "Automatic build"
#    (index 50) and "Share image" (index 51): Reports creation & print,
#    Stations sites, Execution order, Inheritance of item types
#    (based on) + spaces, Production Operator to Work on multiple items
#    simultaneously, Offline work (prod. operator) [with
#    lastRenderedPageBreak], Export to pdf.
# ---------------------------------------------------------------------
$anchor = $d.Paragraphs(50)   # "Automatic build" paragraph

$newParas = @(
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Reports creation &amp; print</w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Stations, sites</w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Execution order</w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Inheritance of item types (based on)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">   </w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Production Operator to Work on multiple items simultaneously</w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Offline work (prod. operator)</w:t></w:r>',
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Export to pdf</w:t></w:r>'
)

foreach ($body in $newParas) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>$body</w:p>"
    $anchor.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 10) "Share image" -> "Share images" (append an "s" run). This is the
#     paragraph right after the 7 freshly inserted ones (originally #51).
# ---------------------------------------------------------------------
$p = $anchor.Next()
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Share image</w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>s</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 11) Insert new "MongoDB instead of google" paragraph (with gramStart/
#     gramEnd proofErr around "google") right after "Share images" and
#     before "AI" (originally #52).
# ---------------------------------------------------------------------
$p.Range.InsertParagraphAfter()
$p = $p.Next()
$xml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">MongoDB instead of </w:t></w:r><w:proofErr w:type=`"gramStart`"/><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>google</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 12) Finally, remove the two originally-empty paragraphs near the top:
#     - the one between "...traceability... single platform." and
#       "Production Operator is required..." (was paragraph 3)
#     - the one right after "...eliminate errors in production." (was
#       paragraph 5, now paragraph 4 once the first is gone)
#     Doing this last means every index used above still refers to the
#     paragraph it meant to when it ran.
# ---------------------------------------------------------------------
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(4).Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
